$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new participants (rows 15 and 16) ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Tsembom Percy"
$ws.Range("C15").Value = 692201677
$ws.Range("D15").Value = 300
$ws.Range("E15").Value = "cash"
$ws.Range("F15").Value = 0
$ws.Range("G15").Formula = "=D15-F15-300"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Tchielong Gaius"
$ws.Range("C16").Value = 695044180
$ws.Range("D16").Value = 1000
$ws.Range("E16").Value = "cash"
$ws.Range("F16").Value = 500
$ws.Range("G16").Formula = "=D16-F16-300"

# --- Recompute / (re)write the "Amount left to reimburse" formulas for the
#     whole block (rows 3-14 previously had blanks or a partial shared
#     formula group; make sure every row down to the new bottom row 16
#     carries the same formula) ---
for ($r = 3; $r -le 16; $r++) {
    $ws.Range("G$r").Formula = "=D$r-F$r-300"
}

# --- Update the active selection to reflect where the user ended up editing ---
$ws.Range("G17").Select()

Write-Host "done"
